$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Confirmados")
$ws2 = $wb.Worksheets.Item("Mortes")

# New date string for row 47 (2020-04-10), added to both data sheets.
# Use the quote-prefix trick so Excel stores it as text (matching the
# existing column A date strings) rather than auto-converting to a date
# serial number, then reset the style to Normal so no extra formatting
# is attached to the cell.
$ws1.Range("A47").Value = "'2020-04-10"
$ws1.Range("A47").Style = "Normal"

$ws2.Range("A47").Value = "'2020-04-10"
$ws2.Range("A47").Style = "Normal"

# Confirmados (sheet1) row 47 values
$ws1.Range("B47").Value = 70
$ws1.Range("C47").Value = 45
$ws1.Range("D47").Value = 166
$ws1.Range("E47").Value = 981
$ws1.Range("F47").Value = 604
$ws1.Range("G47").Value = 1478
$ws1.Range("H47").Value = 555
$ws1.Range("I47").Value = 300
$ws1.Range("J47").Value = 191
$ws1.Range("K47").Value = 293
$ws1.Range("L47").Value = 112
$ws1.Range("M47").Value = 97
$ws1.Range("N47").Value = 698
$ws1.Range("O47").Value = 170
$ws1.Range("P47").Value = 79
$ws1.Range("Q47").Value = 643
$ws1.Range("R47").Value = 684
$ws1.Range("S47").Value = 40
$ws1.Range("T47").Value = 2464
$ws1.Range("U47").Value = 263
$ws1.Range("V47").Value = 636
$ws1.Range("W47").Value = 32
$ws1.Range("X47").Value = 63
$ws1.Range("Y47").Value = 693
$ws1.Range("Z47").Value = 8216
$ws1.Range("AA47").Value = 42
$ws1.Range("AB47").Value = 23

# Mortes (sheet2) row 47 values
$ws2.Range("B47").Value = 2
$ws2.Range("C47").Value = 3
$ws2.Range("D47").Value = 2
$ws2.Range("E47").Value = 50
$ws2.Range("F47").Value = 19
$ws2.Range("G47").Value = 58
$ws2.Range("H47").Value = 14
$ws2.Range("I47").Value = 7
$ws2.Range("J47").Value = 8
$ws2.Range("K47").Value = 16
$ws2.Range("L47").Value = 2
$ws2.Range("M47").Value = 2
$ws2.Range("N47").Value = 17
$ws2.Range("O47").Value = 9
$ws2.Range("P47").Value = 11
$ws2.Range("Q47").Value = 25
$ws2.Range("R47").Value = 65
$ws2.Range("S47").Value = 7
$ws2.Range("T47").Value = 147
$ws2.Range("U47").Value = 11
$ws2.Range("V47").Value = 14
$ws2.Range("W47").Value = 2
$ws2.Range("X47").Value = 3
$ws2.Range("Y47").Value = 18
$ws2.Range("Z47").Value = 540
$ws2.Range("AA47").Value = 4
$ws2.Range("AB47").Value = 0
